$wb = $excel.ActiveWorkbook

# Rename "Table 1" -> "Transship"
$ws = $wb.Worksheets.Item("Table 1")
$ws.Name = "Transship"

# Delete the now-unused empty "Sheet1"
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet1").Delete() | Out-Null
$excel.DisplayAlerts = $true

# Move the selection on the remaining sheet from I13 to G23
$ws.Activate()
$ws.Range("G23").Select() | Out-Null
